$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# --- Row 8: Power indicator (LED) ---
# Mark moves from NOK to OK, comment text tweaked (dropped trailing "so...")
$ws.Range("F8").Value = "Not necessary? Maybe for charging batteries indicatot but the device of the charge already has indications"

# --- Row 10: Components power supplies (dual supply, voltage level) ---
$ws.Range("F10").Value = "Power supply is provided by the msp-fet JTAG_VCC"

# --- Row 15: uC programming/debug connector signals ---
$ws.Range("F15").Value = "Same used for the engeenering models"

# --- Row 16: Decoupling capacitors (10nF/100nF) on power supplies outputs ---
$ws.Range("F16").Value = "Power from debugger must be stable enough"

# --- Row 12: Electrical junctions ---
# Mark moves from NOK to OK, new comment added
$ws.Range("F12").Value = "No errors found when compiling"

# --- Row 19: Differential pair rules ---
$ws.Range("F19").Value = "No high speed signals in design"

# --- Row 21: ESD, EFT, Surge Protection ---
# Mark moves from NOK to OK, new comment added
$ws.Range("F21").Value = "Not needed in this particular board"

# --- Row 23: Add Port Cross Reference (R+P+D) ---
$ws.Range("F23").Value = "Very simple design (1 sheet)"

# --- Row 9: Initial state of control signals (text unchanged, kept for completeness) ---
$ws.Range("F9").Value = "No Integrated Circuits in design"

# --- Row 11: Components polarity (text unchanged, kept for completeness) ---
$ws.Range("F11").Value = "No polarity components"

# --- Row 14: RC circuit (10k - 100nF) o uC RESET pin (text unchanged, kept for completeness) ---
$ws.Range("F14").Value = "No Integrated Circuits in design"

# --- OK/NOK/N-A marks: rows 8, 12 and 21 move from the NOK column (D) to the OK column (C) ---
$ws.Range("C8").Value = "X"
$ws.Range("D8").Value = ""

$ws.Range("C12").Value = "X"
$ws.Range("D12").Value = ""

$ws.Range("C21").Value = "X"
$ws.Range("D21").Value = ""

# --- Update the view: scroll so row 4 is the top visible row, and leave the
#     active selection on F26 (matches the author re-checking the Layout
#     section after finishing the Schematic section). ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F26").Select()
